# Applies the cryptos.xlsx price/volume refresh described in the commit:
# "Updated cryptos list on Sat Aug 17 05:28:07 UTC 2024 with GitHub Actions"
#
# The sheet is a scraped coin-ranking table (columns: A=rank index, B=Coin,
# C=Link, D=Price, E=Volume(1h)). A later scrape produced new Price/Volume
# text values for most rows, and rows 50-51 (RenderToken / VeChain) swapped
# ranking order, so their Coin/Link/Price/Volume values moved between rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds plain-text numeric-looking strings (thousands
# separated with "." like "58.983.69", or plain decimals like "522.32").
# Force the column to Text format first so assigning a decimal-looking
# string (e.g. "522.32", "20.50") keeps its original text/trailing-zero
# formatting instead of being auto-coerced into a real number.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "58.983.69"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("D3").Value = "2.586.76"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "522.32"
$ws.Range("E5").Value = "  +0.79%  "
$ws.Range("D6").Value = "139.32"
$ws.Range("E6").Value = "  -2.27%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "0.565"
$ws.Range("D9").Value = "2.596.72"
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("E10").Value = "  -3.58%  "
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("D13").Value = "0.135"
$ws.Range("E13").Value = "  +3.42%  "
$ws.Range("D14").Value = "3.044.27"
$ws.Range("E14").Value = "  +0.32%  "
$ws.Range("D15").Value = "58.934.20"
$ws.Range("E15").Value = "  +1.58%  "
$ws.Range("D16").Value = "20.50"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("D17").Value = "2.612.92"
$ws.Range("E17").Value = "  +1.23%  "
$ws.Range("E18").Value = "  -0.73%  "
$ws.Range("D19").Value = "338.38"
$ws.Range("E19").Value = "  -1.26%  "
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("E21").Value = "  -1.94%  "
$ws.Range("E22").Value = "  +2.25%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "66.11"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("E25").Value = "  +1.77%  "
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("D28").Value = "7.02"
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  -2.68%  "
$ws.Range("E31").Value = "  -5.33%  "
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("D33").Value = "18.69"
$ws.Range("E33").Value = "  -0.29%  "
$ws.Range("D34").Value = "149.07"
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("E35").Value = "  -0.94%  "
$ws.Range("E36").Value = "  -1.87%  "
$ws.Range("D37").Value = "36.78"
$ws.Range("E37").Value = "  +2.28%  "
$ws.Range("E38").Value = "  +1.30%  "
$ws.Range("D39").Value = "0.827"
$ws.Range("E39").Value = "  -0.84%  "
$ws.Range("D40").Value = "0.817"
$ws.Range("E40").Value = "  -5.87%  "
$ws.Range("E41").Value = "  -0.64%  "
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").Value = "272.25"
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").Value = "10.78"
$ws.Range("E44").Value = "  +1.11%  "
$ws.Range("D45").Value = "0.0954"
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("D47").Value = "0.0517"
$ws.Range("E47").Value = "  -1.24%  "
$ws.Range("D48").Value = "18.40"
$ws.Range("E48").Value = "  -2.33%  "
$ws.Range("D49").Value = "1.962.64"
$ws.Range("E49").Value = "  -0.46%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "0.0220"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "4.50"
$ws.Range("E51").Value = "  -2.04%  "
